$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholder ("datetimeFigureOut" field) on the Slide Master and on
#    every Custom Layout: 09/04/2024 -> 11/04/2024
# ---------------------------------------------------------------------------
$newDate = "11/04/2024"

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq "09/04/2024") {
            $tr.Text = $newDate
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "09/04/2024") {
                $tr.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) "Session 15" / "Session 15 - T" textboxes on every slide ->
#    "Session 16" / "Session 16 - T". Only the trailing run (after
#    the word "Session") is touched so the "Session" run keeps its own
#    formatting untouched.
#    (the COM text getter mangles the non-breaking space to U+FFFD when
#    read back, so detect the runs via StartsWith/EndsWith/Length instead
#    of an exact string match)
# ---------------------------------------------------------------------------
$newTitleSuffix = " 16 - T"
$newFooterSuffix = " 16"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            $full = $tr.Text
            if ($full.StartsWith("Session") -and $full.EndsWith("15 - T") -and $full.Length -eq 14) {
                $sub = $tr.Characters(8, 7)
                $sub.Text = $newTitleSuffix
            } elseif ($full.StartsWith("Session") -and $full.EndsWith("15") -and $full.Length -eq 10) {
                $sub = $tr.Characters(8, 3)
                $sub.Text = $newFooterSuffix
            }
        }
    }
}
